$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (text labels -> shared strings)
$ws.Range("A1").Value = "K"
$ws.Range("C1").Value = "Ln_Prob"

# Data rows (A = K iteration count, C = Ln Prob of data)
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 3).Value = -300472584.39320999
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 3).Value = -300472457.15684301
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 3).Value = -300472452.21204299
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 3).Value = -300472661.58326
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 3).Value = -300472661.58326602
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 3).Value = -300472539.10000002
$ws.Cells.Item(8, 1).Value = 2
$ws.Cells.Item(8, 3).Value = -300472614.89999998
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 3).Value = -300472521.06
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 3).Value = -300472521.06
$ws.Cells.Item(11, 1).Value = 2
$ws.Cells.Item(11, 3).Value = -300472542.69999999
$ws.Cells.Item(12, 1).Value = 3
$ws.Cells.Item(12, 3).Value = -296232690.30000001
$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 3).Value = -296203633.89999998
$ws.Cells.Item(14, 1).Value = 3
$ws.Cells.Item(14, 3).Value = -296203633.89999998
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(15, 3).Value = -296232660.89999998
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 3).Value = -296232641.5
$ws.Cells.Item(17, 1).Value = 3
$ws.Cells.Item(17, 3).Value = -296926334.05000001
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 3).Formula = "=-296310754.1"
$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 3).Value = -296206292.86000001
$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 3).Value = -296206304.04000002
$ws.Cells.Item(21, 1).Value = 3
$ws.Cells.Item(21, 3).Value = -296206718.19999999
$ws.Cells.Item(22, 1).Value = 4
$ws.Cells.Item(22, 3).Value = -292342966.60000002
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 3).Value = -292055931.23671103
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 3).Value = -292704113.94526702
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 3).Value = -292704113.94526702
$ws.Cells.Item(26, 1).Value = 4
$ws.Cells.Item(26, 3).Value = -292273276.63093001
$ws.Cells.Item(27, 1).Value = 4
$ws.Cells.Item(27, 3).Value = -292060816.89999998
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 3).Value = -292060816.89999998
$ws.Cells.Item(29, 1).Value = 4
$ws.Cells.Item(29, 3).Value = -292060816.89999998
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 3).Value = -292689450.80000001
$ws.Cells.Item(31, 1).Value = 4
$ws.Cells.Item(31, 3).Value = -292689607.30000001
$ws.Cells.Item(32, 1).Value = 5
$ws.Cells.Item(32, 3).Value = -288705996.89999998
$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 3).Value = -288439777.38720101
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 3).Value = -288922400.89754599
$ws.Cells.Item(35, 1).Value = 5
$ws.Cells.Item(35, 3).Value = -288922400.89754599
$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 3).Value = -288487401.5
$ws.Cells.Item(37, 1).Value = 5
$ws.Cells.Item(37, 3).Value = -288922400.80000001
$ws.Cells.Item(38, 1).Value = 5
$ws.Cells.Item(38, 3).Value = -288922400.88999999
$ws.Cells.Item(39, 1).Value = 5
$ws.Cells.Item(39, 3).Value = -288487401.5
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 3).Value = -288487401.5
$ws.Cells.Item(41, 1).Value = 5
$ws.Cells.Item(41, 3).Value = -288406818.06
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 3).Value = -285228275.5
$ws.Cells.Item(43, 1).Value = 6
$ws.Cells.Item(43, 3).Value = -285125760.20888197
$ws.Cells.Item(44, 1).Value = 6
$ws.Cells.Item(44, 3).Value = -285296878.96076202
$ws.Cells.Item(45, 1).Value = 6
$ws.Cells.Item(45, 3).Value = -285296878.96076202
$ws.Cells.Item(46, 1).Value = 6
$ws.Cells.Item(46, 3).Value = -285125760.19999999
$ws.Cells.Item(47, 1).Value = 6
$ws.Cells.Item(47, 3).Value = -285296878.95999998
$ws.Cells.Item(48, 1).Value = 6
$ws.Cells.Item(48, 3).Value = -285180088.45999998
$ws.Cells.Item(49, 1).Value = 6
$ws.Cells.Item(49, 3).Value = -285092038.80000001
$ws.Cells.Item(50, 1).Value = 6
$ws.Cells.Item(50, 3).Value = -298509573.19999999
$ws.Cells.Item(51, 1).Value = 6
$ws.Cells.Item(51, 3).Value = -285057013.60000002

# Apply the 0.0 number format to the Ln_Prob column (C1:C51)
$ws.Range("C1:C51").NumberFormat = "0.0"

# Match the recorded selection left by the author
[void]$ws.Range("C2").Select()
